$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Sheet1: VENTAS POR GRUPO
$ws1.Range("M35").Value = 122.42
$ws1.Range("O35").Value = 705.02

$ws1.Range("M60").Value = "2 de 58"
$ws1.Range("N60").Value = "0 de 58"
$ws1.Range("O60").Value = "1 de 58"

# Sheet2: VENTA MENSUAL
$ws2.Range("G2").Value = 3000
$ws2.Range("G3").Value = 2000
$ws2.Range("G5").Value = 500
$ws2.Range("G10").Value = 2000
$ws2.Range("G11").Value = 5000
$ws2.Range("G13").Value = 2000
$ws2.Range("G16").Value = 2500
$ws2.Range("G22").Value = 500
$ws2.Range("G30").Value = 1000
$ws2.Range("G31").Value = 500
$ws2.Range("F35").Value = 827.4400000000001
$ws2.Range("G41").Value = 1000
$ws2.Range("G47").Value = 3000
$ws2.Range("G53").Value = 2500
$ws2.Range("G57").Value = 500
$ws2.Range("G58").Value = 2000
$ws2.Range("F60").Value = 958.98
$ws2.Range("G60").Value = 46000
